# Regenerate s_val data to filter save games.
# Updates columns B, C, D, E, G for rows 2-11 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values keyed by row number: B, C, D, E, G (F is unchanged)
$data = @{
    2  = @(3.286832544864788,     1.655778082260271,    3.537761648806719, 0.4942365360607697, 8.974608811992548)
    3  = @(1.455362044514542,     0.306821227259698,    3.537761648806719, 0.4942365360607697, 5.794181456641729)
    4  = @(0.0006408296065709695, 0.002571899574220771, 0.1494219747398047, 0.4942365360607697, 0.6468712399813661)
    5  = @(3.286832544864788,     1.655778082260271,    3.537761648806719, 0.4942365360607697, 8.974608811992548)
    6  = @(0.01293466051926884,   0.04071648406533734,  3.537761648806719, 0.4942365360607697, 4.085649329452095)
    7  = @(3.286832544864788,     1.655778082260271,    0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    8  = @(1.455362044514542,     0.306821227259698,    0.1494219747398047, 0.4942365360607697, 2.405841782574814)
    9  = @(3.286832544864788,     1.655778082260271,    0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    10 = @(1.455362044514542,     1.655778082260271,    0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    11 = @(0.2917716402565462,    1.655778082260271,    0.7527432677738641, 0.4942365360607697, 3.194529526351451)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G
}
